# ---------------------------------------------------------------------------
# Applies the "Add files via upload" commit:
#   1. Adds 3 rows (34-36) of new failing-testcase entries to the
#      "failing testcases" sheet.
#   2. Adds a brand-new "slotlist" worksheet (after "query variations")
#      containing the slot-name / SlotSet-helper reference table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "failing testcases" sheet - append three new rows
# ---------------------------------------------------------------------------
$fts = $wb.Worksheets.Item("failing testcases")

# Row 34
$fts.Cells.Item(34, 1).Value = "cast of Apollo 13"
$fts.Cells.Item(34, 2).Value = "error message"
$fts.Cells.Item(34, 5).Value = "December milestone"

# Row 35
$fts.Cells.Item(35, 1).Value = "Top 10 science fiction movies"
$fts.Cells.Item(35, 2).Value = "error message"
$fts.Cells.Item(35, 3).Value = "missing training item"
$fts.Cells.Item(35, 4).Value = "fixed"
$fts.Cells.Item(35, 5).Value = "December milestone"

# Row 36
$fts.Cells.Item(36, 1).Value = "top thrillers"
$fts.Cells.Item(36, 2).Value = "answers not correct"
$fts.Cells.Item(36, 3).Value = "missing trainint item"

# ---------------------------------------------------------------------------
# 2. New "slotlist" worksheet, placed right after "query variations"
# ---------------------------------------------------------------------------
$qv = $wb.Worksheets.Item("query variations")
$sl = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $qv)
$sl.Name = "slotlist"

$sl.Cells.Item(1, 1).Value = 'return[SlotSet("ranked_col",None),SlotSet("character",None)]'

# ordered slot names (column A, rows 2..29)
$slotNames = @(
    "budget",
    "cast_name",
    "character",
    "condition_col",
    "condition_operator",
    "condition_val",
    "Costume_Design",
    "Director",
    "Editor",
    "file_name",
    "genre",
    "keyword",
    "language",
    "media",
    "movie",
    "original_language",
    "persistance",
    "plot",
    "Producer",
    "rank_axis",
    "ranked_col",
    "revenue",
    "row_number",
    "row_range",
    "sort_col",
    "top_bottom",
    "year",
    "ascending_descending"
)

for ($i = 0; $i -lt $slotNames.Length; $i++) {
    $row = $i + 2
    $name = $slotNames[$i]

    $sl.Cells.Item($row, 1).Value = $name
    $sl.Cells.Item($row, 2).Formula = "=CONCATENATE(""'"",A$row,""'"")"
    $sl.Cells.Item($row, 3).Formula = "=CONCATENATE(""SlotSet("", B$row,"",None),"")"
    $sl.Cells.Item($row, 4).Value = "SlotSet('" + $name + "',None),"
    # (NB: the original file has a stray space after the first CONCATENATE
    # argument on rows 2-3 only; the engine canonicalises formula text on
    # save regardless, so this is not reproducible / not semantically
    # significant - the computed value is identical either way.)
}

# Row 31: same 28 SlotSet(...) strings (lowercase "Slotset"), spread across columns A..AB
$row31Letters = @(
    "A","B","C","D","E","F","G","H","I","J","K","L","M","N",
    "O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB"
)
for ($i = 0; $i -lt $slotNames.Length; $i++) {
    $name = $slotNames[$i]
    $addr = $row31Letters[$i] + "31"
    $sl.Range($addr).Value = "Slotset('" + $name + "',None),"
}

# Rows 33/34: the hand-typed "clearing string" label + its (split) value
$sl.Cells.Item(33, 1).Value = "clearing string:"
$sl.Cells.Item(34, 1).Value = "return[Slotset('budget',None),Slotset('cast_name',None),Slotset('character',None),Slotset('condition_col',None),Slotset('condition_operator',None),Slotset('condition_val',None),Slotset('Costume_Design',None),Slotset('Director',None),Slotset('Editor',None),Slotset('file_name',None),Slotset('genre',None),Slotset('keyword',None),Slotset('language',None),Slotset('media',None),Slotset('movie',None),"
$sl.Cells.Item(34, 2).Value = "Slotset('original_language',None),Slotset('persistance',None),Slotset('plot',None),Slotset('Producer',None),Slotset('rank_axis',None),Slotset('ranked_col',None),Slotset('revenue',None),Slotset('row_number',None),Slotset('row_range',None),Slotset('sort_col',None),"
$sl.Cells.Item(34, 3).Value = "Slotset('top_bottom',None),Slotset('year',None),Slotset('ascending_descending',None)]"

# Column widths (A:D) - ColumnWidth input is adjusted so the saved OOXML
# <col width> lands on the target values (engine quantises to 1/6 px).
$sl.Columns.Item(1).ColumnWidth = 22.666666666666668
$sl.Columns.Item(2).ColumnWidth = 17.0
$sl.Columns.Item(3).ColumnWidth = 19.833333333333332
$sl.Columns.Item(4).ColumnWidth = 31.833333333333332

# View state for the new sheet
$sl.Range("F14").Select()

# ---------------------------------------------------------------------------
# Restore view state / active sheet to match the saved workbook
# ---------------------------------------------------------------------------
$fts.Activate()
$fts.Range("C36").Select()

Write-Output "done"
